# Add a new row ("release/6.0.2") to the meta-sheet, mirroring the
# formatting of the preceding data row (row 3) which carries no
# explicit cell style (unlike the header rows 1-2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's formatting down to row 4 first so the new row matches
# the existing "no explicit style" data rows instead of inheriting the
# column-level style.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "release/6.0.2"
$ws.Range("B4").Value = "X"
$ws.Range("C4").Value = "X"
$ws.Range("D4").Value = "X"
